$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3.2
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("W2").Value = 7.5
$ws.Range("AL2").Value = 26
$ws.Range("AT2").Value = 2.5
